$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Col1a1"
$ws.Range("C2").Value = "Itga2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 9.108069666666667
$ws.Range("H2").Value = 27.324209
$ws.Range("I2").Value = 0.00155006418458712
$ws.Range("J2").Value = 0.00155006418458712
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 2.672731
$ws.Range("N2").Value = 8.018193
$ws.Range("O2").Value = 0.5408378022089502
$ws.Range("P2").Value = 0.5408378022089503
$ws.Range("Q2").Value = 24.34342014825967
$ws.Range("R2").Value = 219.090781334337
$ws.Range("S2").Value = 0.0008383333068749062
$ws.Range("T2").Value = 0.0008383333068749064

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Col1a1"
$ws.Range("C3").Value = "Itga2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 9.108069666666667
$ws.Range("H3").Value = 27.324209
$ws.Range("I3").Value = 0.00155006418458712
$ws.Range("J3").Value = 0.00155006418458712
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.627877666666667
$ws.Range("N3").Value = 4.883633
$ws.Range("O3").Value = 0.3294075533620981
$ws.Range("P3").Value = 0.3294075533620982
$ws.Range("Q3").Value = 14.82682319681078
$ws.Range("R3").Value = 133.441408771297
$ws.Range("S3").Value = 0.0005106028505990586
$ws.Range("T3").Value = 0.0005106028505990588

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Col1a1"
$ws.Range("C4").Value = "Itga2"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 9.108069666666667
$ws.Range("H4").Value = 27.324209
$ws.Range("I4").Value = 0.00155006418458712
$ws.Range("J4").Value = 0.00155006418458712
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.01852966666666667
$ws.Range("N4").Value = 0.055589
$ws.Range("O4").Value = 0.003749552123152102
$ws.Range("P4").Value = 0.003749552123152104
$ws.Range("Q4").Value = 0.1687694949001111
$ws.Range("R4").Value = 1.518925454101
$ws.Range("S4").Value = 0.000005812046454340666
$ws.Range("T4").Value = 0.000005812046454340669

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Col1a1"
$ws.Range("C5").Value = "Itga2"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 9.108069666666667
$ws.Range("H5").Value = 27.324209
$ws.Range("I5").Value = 0.00155006418458712
$ws.Range("J5").Value = 0.00155006418458712
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.6226963333333334
$ws.Range("N5").Value = 1.868089
$ws.Range("O5").Value = 0.1260050923057995
$ws.Range("P5").Value = 0.1260050923057995
$ws.Range("Q5").Value = 5.67156158517789
$ws.Range("R5").Value = 51.044054266601
$ws.Range("S5").Value = 0.0001953159806588139
$ws.Range("T5").Value = 0.0001953159806588139

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Col1a1"
$ws.Range("C6").Value = "Itga2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 5771.873535333333
$ws.Range("H6").Value = 17315.620606
$ws.Range("I6").Value = 0.9822909543423312
$ws.Range("J6").Value = 0.9822909543423313
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 2.672731
$ws.Range("N6").Value = 8.018193
$ws.Range("O6").Value = 0.5408378022089502
$ws.Range("P6").Value = 0.5408378022089503
$ws.Range("Q6").Value = 15426.665325965
$ws.Range("R6").Value = 138839.987933685
$ws.Range("S6").Value = 0.5312600808762385
$ws.Range("T6").Value = 0.5312600808762388

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Col1a1"
$ws.Range("C7").Value = "Itga2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 5771.873535333333
$ws.Range("H7").Value = 17315.620606
$ws.Range("I7").Value = 0.9822909543423312
$ws.Range("J7").Value = 0.9822909543423313
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.627877666666667
$ws.Range("N7").Value = 4.883633
$ws.Range("O7").Value = 0.3294075533620981
$ws.Range("P7").Value = 0.3294075533620982
$ws.Range("Q7").Value = 9395.90402299351
$ws.Range("R7").Value = 84563.13620694159
$ws.Range("S7").Value = 0.3235740599596277
$ws.Range("T7").Value = 0.3235740599596278

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Col1a1"
$ws.Range("C8").Value = "Itga2"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 5771.873535333333
$ws.Range("H8").Value = 17315.620606
$ws.Range("I8").Value = 0.9822909543423312
$ws.Range("J8").Value = 0.9822909543423313
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.01852966666666667
$ws.Range("N8").Value = 0.055589
$ws.Range("O8").Value = 0.003749552123152102
$ws.Range("P8").Value = 0.003749552123152104
$ws.Range("Q8").Value = 106.9508926518815
$ws.Range("R8").Value = 962.558033866934
$ws.Range("S8").Value = 0.003683151133407393
$ws.Range("T8").Value = 0.003683151133407394

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Col1a1"
$ws.Range("C9").Value = "Itga2"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 5771.873535333333
$ws.Range("H9").Value = 17315.620606
$ws.Range("I9").Value = 0.9822909543423312
$ws.Range("J9").Value = 0.9822909543423313
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.6226963333333334
$ws.Range("N9").Value = 1.868089
$ws.Range("O9").Value = 0.1260050923057995
$ws.Range("P9").Value = 0.1260050923057995
$ws.Range("Q9").Value = 3594.124486915771
$ws.Range("R9").Value = 32347.12038224194
$ws.Range("S9").Value = 0.1237736623730573
$ws.Range("T9").Value = 0.1237736623730574

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Col1a1"
$ws.Range("C10").Value = "Itga2"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.272029666666667
$ws.Range("H10").Value = 3.816089
$ws.Range("I10").Value = 0.0002164813950916887
$ws.Range("J10").Value = 0.0002164813950916887
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 2.672731
$ws.Range("N10").Value = 8.018193
$ws.Range("O10").Value = 0.5408378022089502
$ws.Range("P10").Value = 0.5408378022089503
$ws.Range("Q10").Value = 3.399793123019667
$ws.Range("R10").Value = 30.598138107177
$ws.Range("S10").Value = 0.0001170813219405163
$ws.Range("T10").Value = 0.0001170813219405164

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Col1a1"
$ws.Range("C11").Value = "Itga2"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.272029666666667
$ws.Range("H11").Value = 3.816089
$ws.Range("I11").Value = 0.0002164813950916887
$ws.Range("J11").Value = 0.0002164813950916887
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 1.627877666666667
$ws.Range("N11").Value = 4.883633
$ws.Range("O11").Value = 0.3294075533620981
$ws.Range("P11").Value = 0.3294075533620982
$ws.Range("Q11").Value = 2.070708685704111
$ws.Range("R11").Value = 18.636378171337
$ws.Range("S11").Value = 0.00007131060670556688
$ws.Range("T11").Value = 0.00007131060670556692

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Col1a1"
$ws.Range("C12").Value = "Itga2"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.272029666666667
$ws.Range("H12").Value = 3.816089
$ws.Range("I12").Value = 0.0002164813950916887
$ws.Range("J12").Value = 0.0002164813950916887
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.01852966666666667
$ws.Range("N12").Value = 0.055589
$ws.Range("O12").Value = 0.003749552123152102
$ws.Range("P12").Value = 0.003749552123152104
$ws.Range("Q12").Value = 0.02357028571344444
$ws.Range("R12").Value = 0.212132571421
$ws.Range("S12").Value = 0.0000008117082745889704
$ws.Range("T12").Value = 0.0000008117082745889709

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Col1a1"
$ws.Range("C13").Value = "Itga2"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.272029666666667
$ws.Range("H13").Value = 3.816089
$ws.Range("I13").Value = 0.0002164813950916887
$ws.Range("J13").Value = 0.0002164813950916887
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.6226963333333334
$ws.Range("N13").Value = 1.868089
$ws.Range("O13").Value = 0.1260050923057995
$ws.Range("P13").Value = 0.1260050923057995
$ws.Range("Q13").Value = 0.7920882093245556
$ws.Range("R13").Value = 7.128793883921
$ws.Range("S13").Value = 0.00002727775817101649
$ws.Range("T13").Value = 0.0000272777581710165

# Row 14
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Col1a1"
$ws.Range("C14").Value = "Itga2"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 93.67702500000001
$ws.Range("H14").Value = 281.031075
$ws.Range("I14").Value = 0.01594250007799006
$ws.Range("J14").Value = 0.01594250007799006
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 2.672731
$ws.Range("N14").Value = 8.018193
$ws.Range("O14").Value = 0.5408378022089502
$ws.Range("P14").Value = 0.5408378022089503
$ws.Range("Q14").Value = 250.3734887052751
$ws.Range("R14").Value = 2253.361398347475
$ws.Range("S14").Value = 0.008622306703896162
$ws.Range("T14").Value = 0.008622306703896163

# Row 15
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Col1a1"
$ws.Range("C15").Value = "Itga2"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 93.67702500000001
$ws.Range("H15").Value = 281.031075
$ws.Range("I15").Value = 0.01594250007799006
$ws.Range("J15").Value = 0.01594250007799006
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 1.627877666666667
$ws.Range("N15").Value = 4.883633
$ws.Range("O15").Value = 0.3294075533620981
$ws.Range("P15").Value = 0.3294075533620982
$ws.Range("Q15").Value = 152.494736877275
$ws.Range("R15").Value = 1372.452631895475
$ws.Range("S15").Value = 0.005251579945165764
$ws.Range("T15").Value = 0.005251579945165766

# Row 16
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Col1a1"
$ws.Range("C16").Value = "Itga2"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 93.67702500000001
$ws.Range("H16").Value = 281.031075
$ws.Range("I16").Value = 0.01594250007799006
$ws.Range("J16").Value = 0.01594250007799006
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.01852966666666667
$ws.Range("N16").Value = 0.055589
$ws.Range("O16").Value = 0.003749552123152102
$ws.Range("P16").Value = 0.003749552123152104
$ws.Range("Q16").Value = 1.735804047575
$ws.Range("R16").Value = 15.622236428175
$ws.Range("S16").Value = 0.00005977723501578019
$ws.Range("T16").Value = 0.00005977723501578021

# Row 17
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Col1a1"
$ws.Range("C17").Value = "Itga2"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 93.67702500000001
$ws.Range("H17").Value = 281.031075
$ws.Range("I17").Value = 0.01594250007799006
$ws.Range("J17").Value = 0.01594250007799006
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.6226963333333334
$ws.Range("N17").Value = 1.868089
$ws.Range("O17").Value = 0.1260050923057995
$ws.Range("P17").Value = 0.1260050923057995
$ws.Range("Q17").Value = 58.33233998507502
$ws.Range("R17").Value = 524.9910598656751
$ws.Range("S17").Value = 0.002008836193912354
$ws.Range("T17").Value = 0.002008836193912354
